$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text changes ---
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"

# --- Fix typo in C3 ---
$ws.Range("C3").Value = "casworker-iac-dwp"

# --- Add new row 6 ---
$ws.Range("A6").Value = "BBA9"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Caseworker-divorce"

# --- Formatting: header row bold, size 10 ---
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.Size = 10

# --- Formatting: data rows size 11 ---
$ws.Range("A2:C6").Font.Size = 11

# --- Selection ---
$ws.Range("A1:C6").Select()
